$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview": two new file rows (8cd77d3d.., c723fd20..) land
# ahead of the ".localization-config" row; existing "Ready for
# handoff" rows move to "In Translation" now that handoff is done.
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Delete()

$wsOverview.Range("B2").Value2 = "In Translation"
$wsOverview.Range("C2").Value2 = "In Translation"
$wsOverview.Range("B3").Value2 = "In Translation"
$wsOverview.Range("C3").Value2 = "In Translation"

$wsOverview.Range("A4").Value2 = "8cd77d3d-0189-4bdf-8984-49e473614c01.md"
$wsOverview.Range("B4").Value2 = "Ready for handoff"
$wsOverview.Range("C4").Value2 = "Ready for handoff"

$wsOverview.Range("A5").Value2 = "c723fd20-1ad7-4d11-8239-06a53acfb8aa.md"
$wsOverview.Range("B5").Value2 = "Ready for handoff"
$wsOverview.Range("C5").Value2 = "Ready for handoff"

$wsOverview.Range("A6").Value2 = ".localization-config"
$wsOverview.Range("B6").Value2 = "Not to be localized"
$wsOverview.Range("C6").Value2 = "Not to be localized"

$e2eBase = "https://github.com/OpenLocalizationTest/oltest/blob/45e73766c097b58e853095affa2951f2dd9d3b4a/e2e/"
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/45e73766c097b58e853095affa2951f2dd9d3b4a/.localization-config"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), ($e2eBase + "16debaaa-c511-4974-8d96-ed93ff2c5dc7.md"), "", "", "16debaaa-c511-4974-8d96-ed93ff2c5dc7.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), ($e2eBase + "3e3cf3ab-2f7d-4f75-b392-07ac428c02bb.md"), "", "", "3e3cf3ab-2f7d-4f75-b392-07ac428c02bb.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), ($e2eBase + "8cd77d3d-0189-4bdf-8984-49e473614c01.md"), "", "", "8cd77d3d-0189-4bdf-8984-49e473614c01.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), ($e2eBase + "c723fd20-1ad7-4d11-8239-06a53acfb8aa.md"), "", "", "c723fd20-1ad7-4d11-8239-06a53acfb8aa.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), $configUrl, "", "", ".localization-config")

# ---------------------------------------------------------------
# Sheet "zh-cn": same two new rows, each with its own handoff
# package (.xlf) + a freshly regenerated handoff datetime; the
# ".localization-config" bookkeeping row shifts from row 4 to row 6.
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Hyperlinks.Delete()

$wsZhCn.Range("B2").Value2 = "In Translation"
$wsZhCn.Range("D2").Value2 = "2016-01-17 03:11:56"

$wsZhCn.Range("B3").Value2 = "In Translation"
$wsZhCn.Range("D3").Value2 = "2016-01-17 03:11:56"

$wsZhCn.Range("A4").Value2 = "8cd77d3d-0189-4bdf-8984-49e473614c01.md"
$wsZhCn.Range("B4").Value2 = "Ready for handoff"
$wsZhCn.Range("C4").Value2 = "8cd77d3d-0189-4bdf-8984-49e473614c01.a620b772a59c742d73e85af372e1ead3a467b0f3.zh-cn.xlf"
$wsZhCn.Range("D4").Value2 = "2016-01-17 03:11:56"
$wsZhCn.Range("G4").Value2 = "0001-01-01 00:00:00"
$wsZhCn.Range("H4").Value2 = "Include"

$wsZhCn.Range("A5").Value2 = "c723fd20-1ad7-4d11-8239-06a53acfb8aa.md"
$wsZhCn.Range("B5").Value2 = "Ready for handoff"
$wsZhCn.Range("C5").Value2 = "c723fd20-1ad7-4d11-8239-06a53acfb8aa.854009b7f164ddda0a6cbca4e19c1949073c2209.zh-cn.xlf"
$wsZhCn.Range("D5").Value2 = "2016-01-17 03:11:56"
$wsZhCn.Range("G5").Value2 = "0001-01-01 00:00:00"
$wsZhCn.Range("H5").Value2 = "Include"

$wsZhCn.Range("A6").Value2 = ".localization-config"
$wsZhCn.Range("B6").Value2 = "Not to be localized"
$wsZhCn.Range("D6").Value2 = "0001-01-01 00:00:00"
$wsZhCn.Range("G6").Value2 = "0001-01-01 00:00:00"
$wsZhCn.Range("H6").Value2 = "Ignored"

$wsZhCnHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ff4e98761916f3d43276c0c5f19db7ceb1fdf140/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), ($e2eBase + "16debaaa-c511-4974-8d96-ed93ff2c5dc7.md"), "", "", "16debaaa-c511-4974-8d96-ed93ff2c5dc7.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), ($wsZhCnHandoffBase + "16debaaa-c511-4974-8d96-ed93ff2c5dc7.1db639581a5654d26c25d02710428693a4b1bca9.zh-cn.xlf"), "", "", "16debaaa-c511-4974-8d96-ed93ff2c5dc7.1db639581a5654d26c25d02710428693a4b1bca9.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), ($e2eBase + "3e3cf3ab-2f7d-4f75-b392-07ac428c02bb.md"), "", "", "3e3cf3ab-2f7d-4f75-b392-07ac428c02bb.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C3"), ($wsZhCnHandoffBase + "3e3cf3ab-2f7d-4f75-b392-07ac428c02bb.bd1cd05fa3f857051cf3c2a3b639efecaf99a875.zh-cn.xlf"), "", "", "3e3cf3ab-2f7d-4f75-b392-07ac428c02bb.bd1cd05fa3f857051cf3c2a3b639efecaf99a875.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), ($e2eBase + "8cd77d3d-0189-4bdf-8984-49e473614c01.md"), "", "", "8cd77d3d-0189-4bdf-8984-49e473614c01.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C4"), ($wsZhCnHandoffBase + "8cd77d3d-0189-4bdf-8984-49e473614c01.a620b772a59c742d73e85af372e1ead3a467b0f3.zh-cn.xlf"), "", "", "8cd77d3d-0189-4bdf-8984-49e473614c01.a620b772a59c742d73e85af372e1ead3a467b0f3.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), ($e2eBase + "c723fd20-1ad7-4d11-8239-06a53acfb8aa.md"), "", "", "c723fd20-1ad7-4d11-8239-06a53acfb8aa.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C5"), ($wsZhCnHandoffBase + "c723fd20-1ad7-4d11-8239-06a53acfb8aa.854009b7f164ddda0a6cbca4e19c1949073c2209.zh-cn.xlf"), "", "", "c723fd20-1ad7-4d11-8239-06a53acfb8aa.854009b7f164ddda0a6cbca4e19c1949073c2209.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), $configUrl, "", "", ".localization-config")

# ---------------------------------------------------------------
# Sheet "de-de": same two new rows, each with its own handoff
# package (.xlf) + a freshly regenerated handoff datetime; the
# ".localization-config" bookkeeping row shifts from row 4 to row 6.
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Hyperlinks.Delete()

$wsDeDe.Range("B2").Value2 = "In Translation"
$wsDeDe.Range("D2").Value2 = "2016-01-17 03:12:07"

$wsDeDe.Range("B3").Value2 = "In Translation"
$wsDeDe.Range("D3").Value2 = "2016-01-17 03:12:07"

$wsDeDe.Range("A4").Value2 = "8cd77d3d-0189-4bdf-8984-49e473614c01.md"
$wsDeDe.Range("B4").Value2 = "Ready for handoff"
$wsDeDe.Range("C4").Value2 = "8cd77d3d-0189-4bdf-8984-49e473614c01.a620b772a59c742d73e85af372e1ead3a467b0f3.de-de.xlf"
$wsDeDe.Range("D4").Value2 = "2016-01-17 03:12:07"
$wsDeDe.Range("G4").Value2 = "0001-01-01 00:00:00"
$wsDeDe.Range("H4").Value2 = "Include"

$wsDeDe.Range("A5").Value2 = "c723fd20-1ad7-4d11-8239-06a53acfb8aa.md"
$wsDeDe.Range("B5").Value2 = "Ready for handoff"
$wsDeDe.Range("C5").Value2 = "c723fd20-1ad7-4d11-8239-06a53acfb8aa.854009b7f164ddda0a6cbca4e19c1949073c2209.de-de.xlf"
$wsDeDe.Range("D5").Value2 = "2016-01-17 03:12:07"
$wsDeDe.Range("G5").Value2 = "0001-01-01 00:00:00"
$wsDeDe.Range("H5").Value2 = "Include"

$wsDeDe.Range("A6").Value2 = ".localization-config"
$wsDeDe.Range("B6").Value2 = "Not to be localized"
$wsDeDe.Range("D6").Value2 = "0001-01-01 00:00:00"
$wsDeDe.Range("G6").Value2 = "0001-01-01 00:00:00"
$wsDeDe.Range("H6").Value2 = "Ignored"

$wsDeDeHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/360c7b041c32d0cac26beff67842b300cff0673a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), ($e2eBase + "16debaaa-c511-4974-8d96-ed93ff2c5dc7.md"), "", "", "16debaaa-c511-4974-8d96-ed93ff2c5dc7.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), ($wsDeDeHandoffBase + "16debaaa-c511-4974-8d96-ed93ff2c5dc7.1db639581a5654d26c25d02710428693a4b1bca9.de-de.xlf"), "", "", "16debaaa-c511-4974-8d96-ed93ff2c5dc7.1db639581a5654d26c25d02710428693a4b1bca9.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), ($e2eBase + "3e3cf3ab-2f7d-4f75-b392-07ac428c02bb.md"), "", "", "3e3cf3ab-2f7d-4f75-b392-07ac428c02bb.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C3"), ($wsDeDeHandoffBase + "3e3cf3ab-2f7d-4f75-b392-07ac428c02bb.bd1cd05fa3f857051cf3c2a3b639efecaf99a875.de-de.xlf"), "", "", "3e3cf3ab-2f7d-4f75-b392-07ac428c02bb.bd1cd05fa3f857051cf3c2a3b639efecaf99a875.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), ($e2eBase + "8cd77d3d-0189-4bdf-8984-49e473614c01.md"), "", "", "8cd77d3d-0189-4bdf-8984-49e473614c01.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C4"), ($wsDeDeHandoffBase + "8cd77d3d-0189-4bdf-8984-49e473614c01.a620b772a59c742d73e85af372e1ead3a467b0f3.de-de.xlf"), "", "", "8cd77d3d-0189-4bdf-8984-49e473614c01.a620b772a59c742d73e85af372e1ead3a467b0f3.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), ($e2eBase + "c723fd20-1ad7-4d11-8239-06a53acfb8aa.md"), "", "", "c723fd20-1ad7-4d11-8239-06a53acfb8aa.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C5"), ($wsDeDeHandoffBase + "c723fd20-1ad7-4d11-8239-06a53acfb8aa.854009b7f164ddda0a6cbca4e19c1949073c2209.de-de.xlf"), "", "", "c723fd20-1ad7-4d11-8239-06a53acfb8aa.854009b7f164ddda0a6cbca4e19c1949073c2209.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), $configUrl, "", "", ".localization-config")

